$d = $word.ActiveDocument

# Locate the unique target paragraph ("Dia chi lien lac khi can bao tin: <<...>>")
# via its distinctive Vietnamese label text, so the edit is robust even if
# paragraph indices shift.
$hit = $d.Content
$found = $hit.Find.Execute("Địa chỉ liên lạc khi cần báo tin", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not locate the 'Dia chi lien lac khi can bao tin' label paragraph"
}

$targetPara = $null
$paras = $d.Paragraphs
for ($i = 1; $i -le $paras.Count; $i++) {
    $p = $paras.Item($i)
    if ($p.Range.Start -le $hit.Start -and $p.Range.End -ge $hit.End) {
        $targetPara = $p
    }
}
if ($null -eq $targetPara) {
    throw "Could not resolve the paragraph containing the label"
}

$pRange = $targetPara.Range
$pText = $pRange.Text

$placeholder = "ThiSinh_DienThoai"
$colonIdx = $pText.IndexOf(": <<")
$tokenIdx = $pText.IndexOf($placeholder)
if ($colonIdx -lt 0 -or $tokenIdx -lt 0) {
    throw "Could not locate the ': <<ThiSinh_DienThoai' run sequence in the target paragraph"
}

# Replace from the colon through the end of the paragraph's visible content
# (i.e. up to, but excluding, the trailing paragraph mark). Extending the
# replacement all the way to the last existing run ("'>>'") avoids the host's
# InsertXML reordering runs when a mid-paragraph range is targeted instead.
$replStart = $pRange.Start + $colonIdx
$replEnd = $pRange.End - 1

$target = $d.Range($replStart, $replEnd)
$origTail = $target.Text
if (-not $origTail.EndsWith(">>")) {
    throw "Unexpected paragraph tail while preparing the replacement: $origTail"
}

# Replace the merged "': <<' + 'ThiSinh_DienT' + 'h' + 'oai'" runs with four
# runs: "': '" (space preserved), "'<<'", the new field name
# "'ThiSinh_DCNhanGiayBao'" carrying an explicit black color, and the
# trailing "'>>'" left exactly as it was - matching the target markup.
$xml = @"
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
                <w:sz w:val="26"/>
                <w:szCs w:val="26"/>
              </w:rPr>
              <w:t xml:space="preserve">: </w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
                <w:sz w:val="26"/>
                <w:szCs w:val="26"/>
              </w:rPr>
              <w:t>&lt;&lt;</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
                <w:color w:val="000000"/>
                <w:sz w:val="26"/>
                <w:szCs w:val="26"/>
              </w:rPr>
              <w:t>ThiSinh_DCNhanGiayBao</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
                <w:sz w:val="26"/>
                <w:szCs w:val="26"/>
              </w:rPr>
              <w:t>&gt;&gt;</w:t>
            </w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
"@

$target.InsertXML($xml)
